# Update forests data - 2025-10-15 12:18
#
# 1) The two rows that were on the "New" sheet move down to the bottom of the
#    "Previously added" sheet (rows 180 and 181).
# 2) The "New" sheet is refreshed with five brand-new listings.

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$wsPrev = $wb.Worksheets.Item("Previously added")
$wsNew  = $wb.Worksheets.Item("New")

# Template cell used to clone formatting (font/fill/number format) onto newly
# created cells without Excel registering extra styles/fonts.
$tplRow = 179

function Set-TextCell($cell, [string]$text) {
    # A leading apostrophe forces Excel to treat the value as literal text,
    # which keeps purely numeric looking strings (e.g. cadastre numbers) and
    # empty strings stored as shared-string text instead of being coerced to
    # a number / cleared.
    $cell.Value = "'" + $text
}

function Copy-Format($srcCell, $dstCell) {
    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------------
# 1) Append the two rows moved from "New" onto the end of "Previously added"
# ---------------------------------------------------------------------------

$prevRows = @(
    @{ Row = 180;
       A = "https://www.ss.com/msg/lv/real-estate/wood/daugavpils-and-reg/naujenes-pag/achid.html";
       B = "4 500 €";
       C = "Daugavpils un raj.";
       D = "1 ha.";
       E = "";
       F = 45943.87986111111 },
    @{ Row = 181;
       A = "https://www.ss.com/msg/lv/real-estate/wood/jekabpils-and-reg/vipes-pag/odlhl.html";
       B = "35 000 €";
       C = "Jēkabpils un raj.";
       D = "5 ha.";
       E = "56960010027";
       F = 45944.61388888889 }
)

foreach ($r in $prevRows) {
    $row = $r.Row
    Set-TextCell $wsPrev.Cells.Item($row, 1) $r.A
    Set-TextCell $wsPrev.Cells.Item($row, 2) $r.B
    Set-TextCell $wsPrev.Cells.Item($row, 3) $r.C
    Set-TextCell $wsPrev.Cells.Item($row, 4) $r.D
    Set-TextCell $wsPrev.Cells.Item($row, 5) $r.E
    $wsPrev.Cells.Item($row, 6).Value = $r.F

    $wsPrev.Hyperlinks.Add($wsPrev.Cells.Item($row, 1), $r.A) | Out-Null
}

for ($c = 1; $c -le 6; $c++) {
    Copy-Format $wsPrev.Cells.Item($tplRow, $c) $wsPrev.Cells.Item(180, $c)
    Copy-Format $wsPrev.Cells.Item($tplRow, $c) $wsPrev.Cells.Item(181, $c)
}

# ---------------------------------------------------------------------------
# 2) Replace the contents of "New" with the five freshly scraped listings
# ---------------------------------------------------------------------------

# Drop every existing hyperlink on the sheet (and its now stale relationship)
# before writing new content, so the two old rows don't leave orphaned
# hyperlink/relationship entries behind.
$wsNew.Hyperlinks.Delete() | Out-Null

$newRows = @(
    @{ Row = 2;
       A = "https://www.ss.com/msg/lv/real-estate/wood/aluksne-and-reg/gaujienas-pag/gebgx.html";
       B = "26 900 €";
       C = "Alūksne un raj.";
       D = "4.90 ha.";
       E = "";
       F = 45944.70833333333 },
    @{ Row = 3;
       A = "https://www.ss.com/msg/lv/real-estate/wood/daugavpils-and-reg/salienas-pag/bccgkj.html";
       B = "40 000 €";
       C = "Daugavpils un raj.";
       D = "11.48 ha.";
       E = "44840090024";
       F = 45944.8625 },
    @{ Row = 4;
       A = "https://www.ss.com/msg/lv/real-estate/wood/jekabpils-and-reg/rubenes-pag/indke.html";
       B = "23 500 €";
       C = "Jēkabpils un raj.";
       D = "5 ha.";
       E = "56820070009";
       F = 45945.43611111111 },
    @{ Row = 5;
       A = "https://www.ss.com/msg/lv/real-estate/wood/kraslava-and-reg/kalniesu-pag/iepek.html";
       B = "21 000 €";
       C = "Krāslava un raj.";
       D = "8 ha.";
       E = "60680010191";
       F = 45944.65763888889 },
    @{ Row = 6;
       A = "https://www.ss.com/msg/lv/real-estate/wood/limbadzi-and-reg/ledurgas-pag/djxgk.html";
       B = "49 500 €";
       C = "Limbaži un raj.";
       D = "33 ha.";
       E = "";
       F = 45945.37638888889 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    Set-TextCell $wsNew.Cells.Item($row, 1) $r.A
    Set-TextCell $wsNew.Cells.Item($row, 2) $r.B
    Set-TextCell $wsNew.Cells.Item($row, 3) $r.C
    Set-TextCell $wsNew.Cells.Item($row, 4) $r.D
    Set-TextCell $wsNew.Cells.Item($row, 5) $r.E
    $wsNew.Cells.Item($row, 6).Value = $r.F

    $wsNew.Hyperlinks.Add($wsNew.Cells.Item($row, 1), $r.A) | Out-Null
}

for ($c = 1; $c -le 6; $c++) {
    for ($row = 2; $row -le 6; $row++) {
        Copy-Format $wsPrev.Cells.Item($tplRow, $c) $wsNew.Cells.Item($row, $c)
    }
}

# ---------------------------------------------------------------------------
# Clean up the auto-generated "Hyperlink" named style so the workbook keeps
# using only its original direct-formatting styles.
# ---------------------------------------------------------------------------
foreach ($style in $wb.Styles) {
    if ($style.Name -eq "Hyperlink") {
        $style.Delete() | Out-Null
    }
}

$wb.Save()
